# Applies: Added DegreeRequirement and studenTerm and ApplicationDbContext
# Concretely: builds a "new Degree{...}" C#-snippet helper column on the
# Degree sheet, and restores the active sheet/selection back to Degree.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Degree")

# --- header row (D1:F1) ---------------------------------------------------
$ws.Range("D1").Value = "opening"
$ws.Range("E1").Value = "ending"
$ws.Range("F1").Value = "str"

# --- data rows (D2:D5, E2:E5) ---------------------------------------------
$ws.Range("D2").Value = "new Degree{"
$ws.Range("E2").Value = "},"
$ws.Range("D3").Value = "new Degree{"
$ws.Range("E3").Value = "},"
$ws.Range("D4").Value = "new Degree{"
$ws.Range("E4").Value = "},"
$ws.Range("D5").Value = "new Degree{"
$ws.Range("E5").Value = "},"

# --- formula that stitches the C# line together ----------------------------
$ws.Range("F2").Formula = '=D2&$A$1&"="&A2&E2'

# --- column widths ----------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 33.0
$ws.Columns.Item(3).ColumnWidth = 36.833333333333336
$ws.Columns.Item(4).ColumnWidth = 13.666666666666666
$ws.Columns.Item(5).ColumnWidth = 8.5

# --- header styling (bold header font w/ fill applied, like the rest) ------
$ws.Range("D1:F1").Style = "Normal 2"
$ws.Range("D1:F1").Font.Bold = $true

$ws.Range("E2:E5").Style = "Normal 2"
$ws.Range("F2").Style = "Normal 2"

# --- restore the StudentTerm sheet's own (now inactive) selection ----------
$ws2 = $wb.Worksheets.Item("StudentTerm")
$ws2.Activate()
$ws2.Range("D18").Select()

# --- make Degree the active sheet / selection -------------------------------
$ws.Activate()
$ws.Range("F2:H2").Select()
